$wb = $excel.ActiveWorkbook
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-24 02:33:44"
$zhcn.Range("H3").Value = "2016-03-24 02:34:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-24 02:33:53"
$dede.Range("H3").Value = "2016-03-24 02:35:00"
